# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the per-job worksheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 1320.8572
$ws.Range("I38").Value = 456.16666
$ws.Range("K38").Value = 1368.49998
$ws.Range("M38").Value = -996.4999800000001
# Row 99
$ws.Range("H99").Value = 1416.6666
$ws.Range("I99").Value = 1375
$ws.Range("J99").ClearContents()
$ws.Range("K99").Value = 4125
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = -2627
$ws.Range("N99").Value = -7496
# Row 100
$ws.Range("H100").Value = 1007338.7
$ws.Range("I100").Value = 2503596.8
$ws.Range("J100").Value = 9833.333000000001
$ws.Range("K100").Value = 2503596.8
$ws.Range("L100").Value = 9833.333000000001
$ws.Range("M100").Value = -2503055.8
$ws.Range("N100").Value = -10915.333
# Row 127
$ws.Range("H127").Value = 1999
$ws.Range("I127").Value = 1999
$ws.Range("K127").Value = 5997
$ws.Range("M127").Value = -1037
# Row 129
$ws.Range("H129").Value = 2338.2
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
# Row 132
$ws.Range("H132").Value = 956.9459000000001
$ws.Range("I132").Value = 956.9459000000001
$ws.Range("K132").Value = 2870.8377
$ws.Range("M132").Value = -340.8377
# Row 137
$ws.Range("H137").Value = 2500.7878
$ws.Range("I137").Value = 1488.8948
$ws.Range("K137").Value = 4466.6844
$ws.Range("M137").Value = -1916.6844
# Row 138
$ws.Range("H138").Value = 3946.257
$ws.Range("I138").Value = 3648.8572
$ws.Range("J138").Value = 4144.524
$ws.Range("K138").Value = 10946.5716
$ws.Range("L138").Value = 12433.572
$ws.Range("M138").Value = -5806.571599999999
$ws.Range("N138").Value = -22713.572

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6966.8735
$ws.Range("I32").Value = 6649.9604
$ws.Range("J32").Value = 14995.333
$ws.Range("K32").Value = 6649.9604
$ws.Range("L32").Value = 14995.333
$ws.Range("M32").Value = -6362.9604
$ws.Range("N32").Value = -15569.333
# Row 45
$ws.Range("H45").Value = 2746.5
$ws.Range("I45").Value = 2746.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2746.5
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -2369.5
# Row 110
$ws.Range("H110").Value = 919.875
$ws.Range("I110").Value = 874.8333
$ws.Range("K110").Value = 874.8333
$ws.Range("M110").Value = 1170.1667
# Row 132
$ws.Range("H132").Value = 2810.625
$ws.Range("I132").Value = 2091.1875
$ws.Range("K132").Value = 6273.5625
$ws.Range("M132").Value = -3743.5625

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 925.5333000000001
$ws.Range("I107").Value = 832.5
$ws.Range("J107").Value = 1297.6666
$ws.Range("K107").Value = 832.5
$ws.Range("L107").Value = 1297.6666
$ws.Range("M107").Value = 1087.5
$ws.Range("N107").Value = -5137.6666
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2787.4443
$ws.Range("I31").Value = 2724.7144
$ws.Range("J31").Value = 3007
$ws.Range("K31").Value = 2724.7144
$ws.Range("L31").Value = 3007
$ws.Range("M31").Value = -2429.7144
$ws.Range("N31").Value = -3597
# Row 34
$ws.Range("H34").Value = 2787.4443
$ws.Range("I34").Value = 2724.7144
$ws.Range("J34").Value = 3007
$ws.Range("K34").Value = 2724.7144
$ws.Range("L34").Value = 3007
$ws.Range("M34").Value = -2522.7144
$ws.Range("N34").Value = -3411

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 722.25
$ws.Range("I5").Value = 706.8
$ws.Range("J5").Value = 748
$ws.Range("K5").Value = 2120.4
$ws.Range("L5").Value = 2244
$ws.Range("M5").Value = -2008.4
$ws.Range("N5").Value = -2468
# Row 11
$ws.Range("H11").Value = 2740
$ws.Range("I11").Value = 899.3333
$ws.Range("J11").Value = 5501
$ws.Range("K11").Value = 2697.9999
$ws.Range("L11").Value = 16503
$ws.Range("M11").Value = -2557.9999
$ws.Range("N11").Value = -16783
# Row 18
$ws.Range("H18").Value = 4140
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15338
# Row 23
$ws.Range("H23").Value = 266.33334
$ws.Range("J23").Value = 266.33334
$ws.Range("L23").Value = 799.0000200000001
$ws.Range("N23").Value = -1269.00002
# Row 69
$ws.Range("H69").Value = 1615.3334
$ws.Range("J69").Value = 1999.75
$ws.Range("L69").Value = 5999.25
$ws.Range("N69").Value = -7621.25
# Row 72
$ws.Range("H72").Value = 1615.3334
$ws.Range("J72").Value = 1999.75
$ws.Range("L72").Value = 17997.75
$ws.Range("N72").Value = -26109.75
# Row 75
$ws.Range("H75").Value = 1832.5
$ws.Range("I75").Value = 1300
$ws.Range("J75").Value = 2010
$ws.Range("K75").Value = 3900
$ws.Range("L75").Value = 6030
$ws.Range("M75").Value = -2902
$ws.Range("N75").Value = -8026
# Row 78
$ws.Range("H78").Value = 1832.5
$ws.Range("I78").Value = 1300
$ws.Range("J78").Value = 2010
$ws.Range("K78").Value = 11700
$ws.Range("L78").Value = 18090
$ws.Range("M78").Value = -6708
$ws.Range("N78").Value = -28074
# Row 92
$ws.Range("H92").Value = 975
$ws.Range("J92").Value = 975
$ws.Range("L92").Value = 2925
$ws.Range("N92").Value = -5421
# Row 109
$ws.Range("H109").Value = 2511
$ws.Range("I109").Value = 1277.5
$ws.Range("K109").Value = 3832.5
$ws.Range("M109").Value = -2792.5
# Row 135
$ws.Range("H135").Value = 722.25
$ws.Range("I135").Value = 706.8
$ws.Range("J135").Value = 748
$ws.Range("K135").Value = 6361.2
$ws.Range("L135").Value = 6732
$ws.Range("M135").Value = -3826.2
$ws.Range("N135").Value = -11802
# Row 140
$ws.Range("H140").Value = 1234.3125
$ws.Range("I140").Value = 1234.3125
$ws.Range("K140").Value = 3702.9375
$ws.Range("M140").Value = 1477.0625

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 624.1667
$ws.Range("I97").Value = 572.3333
$ws.Range("J97").Value = 779.6667
$ws.Range("K97").Value = 572.3333
$ws.Range("L97").Value = 779.6667
$ws.Range("M97").Value = -76.33330000000001
$ws.Range("N97").Value = -1771.6667
# Row 113
$ws.Range("H113").Value = 2152.6667
$ws.Range("I113").Value = 1480.5
$ws.Range("K113").Value = 1480.5
$ws.Range("M113").Value = 689.5
# Row 114
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 557.6667
$ws.Range("I107").Value = 509.2
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1527.6
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 392.4000000000001
$ws.Range("N107").Value = -6240
# Row 122
$ws.Range("H122").Value = 4666.3335
$ws.Range("I122").Value = 4666.3335
$ws.Range("K122").Value = 13999.0005
$ws.Range("M122").Value = -11549.0005
